# Capitalize the first letter of each TYPE value in the "Field List" sheet
# (column D, rows 2-66): "number" -> "Number", "text" -> "Text", "date" -> "Date"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Field List")

for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $value = $cell.Value2
    if ($value -ne $null -and $value -ne "") {
        $newValue = $value.Substring(0,1).ToUpper() + $value.Substring(1)
        $cell.Value2 = $newValue
    }
}
